$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-15 (columns D, L, M, N, O, P, Q, S, T).
# All other columns (A,B,C,E,F,G,H,I,J,K,R) are unchanged across rows.
$rows = @{
  2  = @{ D = 44904; L = "Primera"; M = 45;  N = 15000; O = 15000; P = 15000; Q = '$/bandeja 10 kilos';            S = 1500; T = 10 }
  3  = @{ D = 44904; L = "Segunda"; M = 60;  N = 10000; O = 10000; P = 10000; Q = '$/bandeja 10 kilos';            S = 1000; T = 10 }
  4  = @{ D = 44391; L = "Primera"; M = 15;  N = 1500;  O = 1500;  P = 1500;  Q = '$/kilo (en caja de 15 kilos)'; S = 1500; T = 1  }
  5  = @{ D = 44391; L = "Segunda"; M = 20;  N = 1000;  O = 1000;  P = 1000;  Q = '$/kilo (en caja de 15 kilos)'; S = 1000; T = 1  }
  6  = @{ D = 44195; L = "Primera"; M = 20;  N = 15000; O = 15000; P = 15000; Q = '$/bandeja 10 kilos';            S = 1500; T = 10 }
  7  = @{ D = 44400; L = "Primera"; M = 25;  N = 1500;  O = 1500;  P = 1500;  Q = '$/kilo (en caja de 15 kilos)'; S = 1500; T = 1  }
  8  = @{ D = 44336; L = "Primera"; M = 10;  N = 1500;  O = 1500;  P = 1500;  Q = '$/kilo (en caja de 15 kilos)'; S = 1500; T = 1  }
  9  = @{ D = 44880; L = "Primera"; M = 200; N = 20000; O = 20000; P = 20000; Q = '$/bandeja 10 kilos';            S = 2000; T = 10 }
  10 = @{ D = 44880; L = "Segunda"; M = 180; N = 15000; O = 15000; P = 15000; Q = '$/bandeja 10 kilos';            S = 1500; T = 10 }
  11 = @{ D = 44343; L = "Primera"; M = 20;  N = 1700;  O = 1700;  P = 1700;  Q = '$/kilo (en caja de 15 kilos)'; S = 1700; T = 1  }
  12 = @{ D = 44292; L = "Primera"; M = 50;  N = 14000; O = 14000; P = 14000; Q = '$/bandeja 10 kilos';            S = 1400; T = 10 }
  13 = @{ D = 44309; L = "Primera"; M = 10;  N = 1600;  O = 1600;  P = 1600;  Q = '$/kilo (en caja de 15 kilos)'; S = 1600; T = 1  }
  14 = @{ D = 44371; L = "Primera"; M = 20;  N = 1800;  O = 1800;  P = 1800;  Q = '$/kilo (en caja de 15 kilos)'; S = 1800; T = 1  }
  15 = @{ D = 44371; L = "Segunda"; M = 30;  N = 1200;  O = 1200;  P = 1200;  Q = '$/kilo (en caja de 15 kilos)'; S = 1200; T = 1  }
}

foreach ($r in $rows.Keys) {
  $v = $rows[$r]
  $ws.Cells.Item($r, 4).Value2  = $v.D   # D: Fecha
  $ws.Cells.Item($r, 12).Value  = $v.L   # L: Calidad
  $ws.Cells.Item($r, 13).Value2 = $v.M   # M: Volumen
  $ws.Cells.Item($r, 14).Value2 = $v.N   # N: Precio minimo
  $ws.Cells.Item($r, 15).Value2 = $v.O   # O: Precio maximo
  $ws.Cells.Item($r, 16).Value2 = $v.P   # P: Precio promedio ponderado
  $ws.Cells.Item($r, 17).Value  = $v.Q   # Q: Unidad de comercializacion
  $ws.Cells.Item($r, 19).Value2 = $v.S   # S: Precio $/Kg
  $ws.Cells.Item($r, 20).Value2 = $v.T   # T: Kg / unidad
}
